$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 1479
$ws.Range("K3").Value = 1408
$ws.Range("K4").Value = 306
$ws.Range("K5").Value = 96
$ws.Range("K6").Value = 1810
$ws.Range("K7").Value = 5099

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 96
$ws.Range("K3").Value = 90
$ws.Range("K6").Value = 101
$ws.Range("K7").Value = 311

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 45
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 57
$ws.Range("K6").Value = 64
$ws.Range("K7").Value = 174

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 35
$ws.Range("K7").Value = 145
$ws.Range("K8").Value = 311
$ws.Range("K11").Value = 102
$ws.Range("K19").Value = 132
$ws.Range("K20").Value = 117
$ws.Range("K23").Value = 52
$ws.Range("K27").Value = 61
$ws.Range("K29").Value = 233
$ws.Range("K31").Value = 59
$ws.Range("K34").Value = 33
$ws.Range("K36").Value = 57
$ws.Range("K37").Value = 174
$ws.Range("K42").Value = 176
$ws.Range("K43").Value = 49
$ws.Range("K44").Value = 47
$ws.Range("K47").Value = 37
$ws.Range("K48").Value = 56
$ws.Range("K51").Value = 59
$ws.Range("K52").Value = 138
$ws.Range("K53").Value = 79
$ws.Range("K63").Value = 19
$ws.Range("K64").Value = 32
$ws.Range("K66").Value = 21
$ws.Range("K67").Value = 195
$ws.Range("K70").Value = 11
$ws.Range("K72").Value = 22
$ws.Range("K73").Value = 53
$ws.Range("K76").Value = 70
$ws.Range("K79").Value = 136
$ws.Range("K83").Value = 103
$ws.Range("K85").Value = 264
$ws.Range("K89").Value = 66
$ws.Range("K95").Value = 88
$ws.Range("K98").Value = 34
$ws.Range("K99").Value = 93
$ws.Range("K101").Value = 5099

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 61
$ws.Range("K7").Value = 195

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 62
$ws.Range("K3").Value = 75
$ws.Range("K7").Value = 233

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 42
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K3").Value = 17
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 70

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 41
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 20
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K3").Value = 47
$ws.Range("K7").Value = 136

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 36
$ws.Range("K7").Value = 117

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 22
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 46
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K5").Value = 2
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 13
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K3").Value = 4
$ws.Range("K6").Value = 21

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 14
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("K5").Value = 3
$ws.Range("K6").Value = 11

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 49

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 97
$ws.Range("K3").Value = 82
$ws.Range("K4").Value = 14
$ws.Range("K6").Value = 67
$ws.Range("K7").Value = 264

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K2").Value = 3
$ws.Range("K3").Value = 4
$ws.Range("K6").Value = 22

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 138
